# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (header is row 1, data starts row 2).
$lastRow = $ws.UsedRange.Rows.Count

# New header cells: AD1="Wins", AE1="Losses", AF1="Ties".
# Copy the formatting from the neighboring header cell (AC1, which uses the
# bold/bordered/centered header style) onto the three new header cells so
# they visually match the rest of row 1, then set their text.
$headerSrc = $ws.Cells.Item(1, 29)
$headerSrc.Copy($ws.Cells.Item(1, 30))
$headerSrc.Copy($ws.Cells.Item(1, 31))
$headerSrc.Copy($ws.Cells.Item(1, 32))

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Fill the team record for every data row (2-49): Wins=74, Losses=87, Ties=0.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
